$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three data values in row 1 (A1, B1, C1)
$ws.Range("A1").Value = 149.06717768846821
$ws.Range("B1").Value = 4.5462262155979705
$ws.Range("C1").Value = 1.163894324853229

# Column A (width 11.7109375) is unchanged; columns B and C should match it.
# Set columns B and C to the same ColumnWidth as column A so their stored
# OOXML width matches column A's.
$ws.Columns.Item(2).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 10.833333333333334
